# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows for "Cebolla" (onion) at Vega Central Mapocho de Santiago,
# just before the existing row 1119, shifting the rest of the data down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 1119 (existing rows 1119:1161 move to 1123:1165)
$ws.Rows.Item(1119).Resize(4).Insert() | Out-Null

# New row data: [Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, Origen, PrecioKg, KgUnidades]
$newRows = @(
    @{ Row = 1119; H = "Sin especificar"; I = "1a (guarda)";  J = 430;  K = 4900; L = 5200; M = 5050; N = "`$/malla 16 kilos"; O = "Región de O'Higgins"; P = 316; Q = 16 },
    @{ Row = 1120; H = "Sin especificar"; I = "1a nueva(o)";  J = 8600; K = 2000; L = 2200; M = 2100; N = "`$/paquete 20 unidades (volumen en unidades)"; O = "Región de O'Higgins"; P = 105; Q = 20 },
    @{ Row = 1121; H = "Sin especificar"; I = "2a (guarda)";  J = 250;  K = 4400; L = 4700; M = 4550; N = "`$/malla 16 kilos"; O = "Región de O'Higgins"; P = 284; Q = 16 },
    @{ Row = 1122; H = "Sin especificar"; I = "2a nueva(o)";  J = 3200; K = 1600; L = 1800; M = 1700; N = "`$/paquete 20 unidades (volumen en unidades)"; O = "Región de O'Higgins"; P = 85;  Q = 20 }
)

$newDate = Get-Date -Year 2021 -Month 11 -Day 9 -Hour 0 -Minute 0 -Second 0

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = $newDate
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112004
    $ws.Cells.Item($row, 7).Value = "Cebolla"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
